$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Remove existing hyperlinks on this sheet before rewriting the data range
$ws.Hyperlinks.Delete()

$data = @(
    ,@(2, '2025-11-14 12:36:27', '建設・土木業界向け施工機械のAI自動制御・アタッチメント開発を支援してくださるエンジニア募集', 'システム開発', '200,000 円 ~ 300,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5434128', 368, '🔥AI,Ai ◆開発')
    ,@(3, '2025-11-14 12:36:27', '企業のMicrosoft Copilot導入・活用支援AIコンサルタント募集(研修講師・メンター)', 'システム開発', '200,000 円 ~ 300,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5434363', 348, '🔥AI,Ai ◆コンサル')
    ,@(4, '2025-11-14 12:36:27', '【GAS】Yahoo!ショッピング注文完了メール (Gmail) からスプレッドシートに転記する仕事', 'システム開発', '20,000 円 ~ 50,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5433649', 298, '🔥AI,Ai')
    ,@(5, '2025-11-14 12:36:27', '画像処理システム(ツール)の開発', 'システム開発', '20,000 円 ~ 50,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5434134', 138, '◆ツール,開発')
    ,@(6, '2025-11-14 12:36:27', '英語教育の公式LINEアカウント開発・運用スタッフ募集【即日〜3月/4ヶ月/継続可能】', 'システム開発', '20,000 円 ~ 50,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5433668', 63, '◆開発')
    ,@(7, '2025-11-14 12:36:27', '【急募】WordPressでの商品検索サイト構築依頼', 'システム開発', '200,000 円 ~ 300,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5433985', 58, '◇サイト ○WordPress')
    ,@(8, '2025-11-14 12:36:27', '【GAS活用】業務改善システムの構築依頼', 'システム開発', '50,000 円 ~ 100,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5434156', 53, '◇業務改善')
    ,@(9, '2025-11-14 12:36:27', 'wordpressレンダリングを妨げるリソースの除外', 'システム開発', '200,000 円 ~ 300,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5016989', 33, '○WordPress')
    ,@(10, '2025-11-14 12:36:27', '【相談から】Laravel7からLaravel12へのサーバーアップデート依頼', 'システム開発', '500,000 円 ~ 1,000,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5433727', 25, $null)
    ,@(11, '2025-11-14 12:36:27', '初回 ★社内の音響設計スキル向上のため、Modeler / EASE Focus を教えていただける方', 'システム開発', '50,000 円 ~ 100,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5433823', 18, $null)
    ,@(12, '2025-11-14 12:36:27', '月1~5万円以内の小規模タスク依頼', 'システム開発', '20,000 円 ~ 50,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5433937', 13, $null)
    ,@(13, '2025-11-14 12:36:27', '【相談のみ】Unityで自動ルート設計プログラムが実現可能か専門家に相談がしたい', 'システム開発', '1,000 ~ 5,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5434061', 10, $null)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
    $ws.Range("F$r").Value = $row[6]
    $ws.Hyperlinks.Add($ws.Range("F$r"), $row[6])
    $ws.Range("F$r").Style = "Hyperlink"
    $ws.Range("G$r").Value = $row[7]
    if ($row[8] -ne $null) {
        $ws.Range("H$r").Value = $row[8]
    } else {
        $ws.Range("H$r").ClearContents()
    }
}
